# "Add files via upload" — re-upload of questions.xlsx with the "Общие
# знания" (General knowledge) section renamed to "Общие_знания" (space
# replaced with underscore) for every question row in that section
# (rows 52-61, column A / "Раздел"), plus the sheet's last-used selection
# moving down to the new last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 52-61 are the "Общие знания" quiz questions (column A = "Раздел").
# Re-writing the value (with the underscore instead of a space) creates a
# new shared-string entry; once every old reference has been replaced the
# now-unused original string is dropped from the table on save, so the
# whole shared-strings table renumbers exactly like the upstream edit.
for ($r = 52; $r -le 61; $r++) {
    $ws.Cells.Item($r, 1).Value = "Общие_знания"
}

# Reflect the author's final on-screen selection/scroll position (the
# last row of the sheet, cell B62) when the file was re-saved.
$ws.Activate()
$ws.Range("B62").Select()
